$wb = $excel.ActiveWorkbook

# ------------------------------------------------------------------
# 1. Nexus sheet: add two new rows (account/password block) and a
#    stray formatted-but-empty cell (B14) that carries a wrap-text
#    style (this becomes cellXfs index 12 = wrapText only).
# ------------------------------------------------------------------
$nexus = $wb.Worksheets.Item("Nexus")
$nexus.Range("A5").Value = "account: "
$nexus.Range("B5").Value = "devops"
$nexus.Range("A6").Value = "password: "
$nexus.Range("B6").Value = "FX9X8el*KLCo!JcB"
$nexus.Range("B14").WrapText = $true

# ------------------------------------------------------------------
# 2. Add the "Jenkins" sheet right after "Nexus".
#    A1 holds the multi-line jenkins credential block, wrapped and
#    top-aligned (cellXfs index 13 = wrapText + vertical top).
# ------------------------------------------------------------------
$lastSheet = $wb.Worksheets.Item($wb.Worksheets.Count)
$jenkins = $wb.Worksheets.Add($null, $lastSheet)
$jenkins.Name = "Jenkins"
$jenkins.Columns.Item(1).ColumnWidth = 62.4140625
$jenkins.Range("A1").Value = "jenkins:  https://jenkins.nissan.com.cn/`naccount: devops`npassword: 8CnFAb@2xgL9O2&*"
$jenkins.Range("A1").WrapText = $true
$jenkins.Range("A1").VerticalAlignment = -4160
$jenkins.Range("A1").RowHeight = 42

# ------------------------------------------------------------------
# 3. Add the "aliyun" sheet right after "Jenkins".
# ------------------------------------------------------------------
$aliyun = $wb.Worksheets.Add($null, $jenkins)
$aliyun.Name = "aliyun"
$aliyun.Columns.Item(1).ColumnWidth = 8.6640625
$aliyun.Columns.Item(2).ColumnWidth = 57.83203125

$aliyun.Range("A1").Value = "url"
$aliyun.Range("B1").Value = "https://signin.aliyun.com/1615541751802351.onaliyun.com/login.htm"
$aliyun.Range("B1").Style = "Hyperlink"

$aliyun.Range("A2").Value = "account"
$aliyun.Range("B2").Value = "devops@1615541751802351.onaliyun.com"
$aliyun.Range("B2").Style = "Hyperlink"

$aliyun.Range("A3").Value = "password"
$aliyun.Range("B3").Value = "4ZH7Pkp`$9Zaa#Pp5"

# ------------------------------------------------------------------
# 4. Restore / update the view selections on every sheet.
# ------------------------------------------------------------------
$ws1 = $wb.Worksheets.Item("40.73.1.192")
$ws1.Range("B19").Select()

$ws2 = $wb.Worksheets.Item("40.73.23.194")
$ws2.Range("F3").Select()

$ws3 = $wb.Worksheets.Item("139.217.228.205")
$ws3.Rows.Item(8).Select()

$nexus.Range("H15:H18").Select()

$jenkins.Range("H10").Select()

$aliyun.Range("B15").Select()

# ------------------------------------------------------------------
# 5. Jenkins ends up the active / selected tab (activeTab=4,
#    tabSelected on sheet5.xml).
# ------------------------------------------------------------------
$jenkins.Activate()
